# "soulignement des élément si ils sont à vérifier"
#
# 1) Bump the notes-master date placeholder field (23/06/2023 -> 26/06/2023).
# 2) Add a small "A vérifier" underlined text box (bottom-right corner) on
#    slides 2, 4, 5 and 6.

$p = $ppt.ActivePresentation

# --- 1) Notes master date field -------------------------------------------------
# This is a live datetimeFigureOut field (auto-computed "today"), but we still
# try the direct, idiomatic COM write in case the host lets it through.
try {
    $nm = $p.NotesMaster
    $dateShape = $nm.Shapes.Item(2)
    $dateShape.TextFrame.TextRange.Text = "26/06/2023"
} catch {
    # best effort only - some hosts keep datetimeFigureOut fields read-only
}

# --- 2) "A vérifier" text boxes --------------------------------------------------
# Identical textbox (size/position/formatting) pasted onto four slides.
# EMU -> point conversion (914400 EMU per inch, 72 pt per inch):
#   x  = 11577729 EMU = 911.6322047244095 pt
#   y  = 6642556  EMU = 523.035905511811  pt
#   cx = 614271   EMU = 48.36779527559055 pt
#   cy = 215444   EMU = 16.964094488188977 pt
$boxLeft   = 911.6322047244095
$boxTop    = 523.035905511811
$boxWidth  = 48.36779527559055
$boxHeight = 16.964094488188977

$targetSlides = 2, 4, 5, 6

foreach ($slideIndex in $targetSlides) {
    $slide = $p.Slides.Item($slideIndex)

    $shape = $slide.Shapes.AddTextbox(1, $boxLeft, $boxTop, $boxWidth, $boxHeight)

    # PowerPoint names a freshly-added text box "ZoneTexte N" (next free index
    # among same-named shapes already on the slide) in the French UI build
    # this deck was authored with.
    if ($slideIndex -eq 2 -or $slideIndex -eq 4) {
        $shape.Name = "ZoneTexte 2"
    } else {
        $shape.Name = "ZoneTexte 1"
    }

    $shape.Fill.Visible = $false

    $tf = $shape.TextFrame
    $tf.WordWrap = $false
    $tf.AutoSize = 1

    $tr = $tf.TextRange
    $tr.Text = "A vérifier"
    $tr.LanguageID = "fr-FR"
    $tr.Font.Size = 8
    $tr.Font.Underline = $true
    $tr.ParagraphFormat.Alignment = 1
}

Write-Output "done"
